$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Groupe" column values: students moved from semester 2 groups to semester 1 groups ---
# (shared strings "2-A"/"2-B"/"2-C" become "1-A"/"1-B"/"1-C" for every student using them)
$groupMap = @{
    "2-A" = "1-A"
    "2-B" = "1-B"
    "2-C" = "1-C"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -and $groupMap.ContainsKey($val)) {
        $cell.Value = $groupMap[$val]
    }
}

# --- 2. Update "Numero" (A) and "Moyenne de l'etudiant" (E) columns per row ---
# Numero: 2015xxxx -> 2017xxxx (registration year correction)
# Moyenne: refreshed average values following the fix
$rowData = @{
    3 = @(20170926, 14)
    4 = @(20170927, 10)
    5 = @(20170928, 9)
    6 = @(20170929, 15)
    7 = @(20170930, 20)
    8 = @(20170931, 19)
    9 = @(20170932, 7)
    10 = @(20170933, 9)
    11 = @(20170934, 20)
    12 = @(20170935, 5)
    13 = @(20170936, 20)
    14 = @(20170937, 16)
    15 = @(20170938, 18)
    16 = @(20170939, 17)
    17 = @(20170940, 13)
    18 = @(20170941, 13)
    19 = @(20170942, 8)
    20 = @(20170943, 16)
    21 = @(20170944, 16)
    22 = @(20170945, 11)
    23 = @(20170946, 13)
    24 = @(20170947, 13)
    25 = @(20170948, 7)
    26 = @(20170949, 10)
    27 = @(20170950, 7)
    28 = @(20170951, 9)
    29 = @(20170952, 15)
    30 = @(20170953, 12)
    31 = @(20170954, 12)
    32 = @(20170955, 16)
    33 = @(20170956, 16)
    34 = @(20170957, 14)
    35 = @(20170958, 12)
    36 = @(20170959, 17)
    37 = @(20170960, 19)
    38 = @(20170961, 9)
    39 = @(20170962, 13)
    40 = @(20170963, 8)
    41 = @(20170964, 16)
    42 = @(20170965, 7)
    43 = @(20170966, 12)
    44 = @(20170967, 12)
    45 = @(20170968, 7)
    46 = @(20170969, 19)
    47 = @(20170970, 10)
    48 = @(20170971, 11)
    49 = @(20170972, 9)
    50 = @(20170973, 10)
    51 = @(20170974, 6)
    52 = @(20170975, 19)
    53 = @(20170976, 10)
    54 = @(20170977, 11)
    55 = @(20170978, 16)
    56 = @(20170979, 6)
    57 = @(20170980, 17)
    58 = @(20170981, 5)
    59 = @(20170982, 6)
    60 = @(20170983, 11)
    61 = @(20170984, 9)
    62 = @(20170985, 18)
    63 = @(20170986, 12)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item([int]$r, 1).Value = $vals[0]
    $ws.Cells.Item([int]$r, 5).Value = $vals[1]
}
